$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Swap country labels (text only; underlying row data stays positioned) ---
$ws.Range("A6").Value  = "Brasil"
$ws.Range("A7").Value  = "España"
$ws.Range("A95").Value = "Gabon"
$ws.Range("A96").Value = "Nueva Zelanda"

# --- Update timestamp message ---
$ws.Range("A1").Value = "Datos actualizados a 21 de Mayo de 2020 a las 00:05"

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 1588683
$ws.Range("C4").Value = 18100
$ws.Range("D4").Value = 362161
$ws.Range("E4").Value = 1131674
$ws.Range("G4").Value = 1315
$ws.Range("H4").Value = 94848

# --- Row 6: now Brasil ---
$ws.Range("B6").Value = 285418
$ws.Range("C6").Value = 13533
$ws.Range("D6").Value = 106794
$ws.Range("E6").Value = 159969
$ws.Range("G6").Value = 672
$ws.Range("H6").Value = 18655

# --- Row 7: now España ---
$ws.Range("B7").Value = 279524
$ws.Range("C7").Value = 721
$ws.Range("D7").Value = 196958
$ws.Range("E7").Value = 54678
$ws.Range("G7").Value = 110
$ws.Range("H7").Value = 27888

# --- Row 11: Alemania ---
$ws.Range("B11").Value = 178494
$ws.Range("C11").Value = 667
$ws.Range("E11").Value = 13329

# --- Row 95: now Gabon ---
$ws.Range("B95").Value = 1567
$ws.Range("C95").Value = 65
$ws.Range("D95").Value = 365
$ws.Range("E95").Value = 1190
$ws.Range("H95").Value = 12

# --- Row 96: now Nueva Zelanda ---
$ws.Range("B96").Value = 1503
$ws.Range("D96").Value = 1447
$ws.Range("E96").Value = 35
$ws.Range("H96").Value = 21

# --- Row 165 ---
$ws.Range("D165").Value = 80
$ws.Range("E165").Value = 36
